$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the rows being edited (rows 8:9), matching the author's selection
# captured in the saved file (activeCell A8, sqref A8:XFD9).
$ws.Range("A8:XFD9").Select()

# Row 8 (Magnesium chloride unit price): update lower bound and
# overwrite the formula-derived Lower/Upper cells with plain values.
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Row 9 (Zinc sulfate unit price): same treatment.
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931
